$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new rows of wallet-tracing data (date + USD value).
#
# The "Date" column stores plain text that merely looks like a date
# ("YYYY-MM-DD"), so before writing those strings we mark the destination
# cells as Text ("@") to stop Excel from auto-converting them into date
# serial numbers. We then restore the cells to the default "Normal" style
# so no stray cell-level formatting is left behind, matching a plain
# shared-string cell with no explicit style.
$ws.Range("A31:A32").NumberFormat = "@"

$ws.Cells.Item(31, 1).Value = "2024-10-01"
$ws.Cells.Item(31, 2).Value = 0.00000043

$ws.Cells.Item(32, 1).Value = "2024-10-05"
$ws.Cells.Item(32, 2).Value = 0.00000041

$ws.Range("A31:A32").Style = "Normal"
